# Updated the RF model for the prediction of sedimentation rates
# The HD_diff column (old column P, = O-N shared formula) is no longer
# needed, so it is deleted entirely. This shifts the old ksed (1/h)
# column (Q) one place to the left, into the new column P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "HD_diff" column (old column P). Deleting the
# entire column shifts every column to its right (just "ksed (1/h)",
# formerly column Q) one place left, preserving each cell's value and
# number format, and drops the now-stale shared formula along with it.
$ws.Range("P1").EntireColumn.Delete()

# Restore the view state left by the editing session: a smaller zoom
# level and the last selected cell on the sheet.
$excel.ActiveWindow.Zoom = 130
$ws.Range("O11").Select()
